$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the unused user guide references ("#201" and "#102") from D5/D6
# while preserving their existing fill/formatting.
$ws.Range("D5").ClearContents()
$ws.Range("D6").ClearContents()

# Update the selected/active cell to reflect where the author left off.
$ws.Range("J9").Select()
